# B1--and-B2-PowerPoint.pptx edit
#
# 1) The single table on slide 5 gets a different (built-in) table style
#    applied - tableStyleId changes from {28E2AF44-098C-4902-8719-5127C7976A43}
#    to {E7F6E40B-10DD-4FA3-A43A-45A9AA0BE6B4}. Table.Style is a read-only
#    reflection of the style id in this host - it has to be changed with
#    Table.ApplyStyle(id).
#
# 2) The deck's design/theme colour palette is switched from the "Integral"
#    (Red Violet) palette back to the stock "Office" palette. The only
#    theme reachable from the object model is the one wired to the slide
#    master (ppt/theme/theme2.xml in the underlying package), so the swap
#    is performed by writing the twelve "Office" theme colours into the
#    master's ThemeColorScheme, slot by slot. .RGB takes/returns a
#    VBA-style 0xBBGGRR value, so each target 0xRRGGBB colour below is
#    byte-swapped before assignment.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{E7F6E40B-10DD-4FA3-A43A-45A9AA0BE6B4}")

# --- 2. Theme colour scheme: Integral/"Red Violet" -> stock "Office" ---
$master = $p.Designs.Item(1).SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Colors(1).RGB  = 0x000000   # dk1
$colors.Colors(2).RGB  = 0xFFFFFF   # lt1
$colors.Colors(3).RGB  = 0x6A5444   # dk2      (0x44546A)
$colors.Colors(4).RGB  = 0xE6E6E7   # lt2      (0xE7E6E6)
$colors.Colors(5).RGB  = 0xD59B5B   # accent1  (0x5B9BD5)
$colors.Colors(6).RGB  = 0x317DED   # accent2  (0xED7D31)
$colors.Colors(7).RGB  = 0xA5A5A5   # accent3  (0xA5A5A5)
$colors.Colors(8).RGB  = 0x00C0FF   # accent4  (0xFFC000)
$colors.Colors(9).RGB  = 0xC47244   # accent5  (0x4472C4)
$colors.Colors(10).RGB = 0x47AD70   # accent6  (0x70AD47)
$colors.Colors(11).RGB = 0xC16305   # hlink    (0x0563C1)
$colors.Colors(12).RGB = 0x724F95   # folHlink (0x954F72)
